# Update the DFII10 series with the latest FRED observation(s):
#  - row 122, col A (date) is corrected from 2025-08-29 (45898) to 2025-09-09 (45909)
#  - a new row 123 is appended: 2025-09-15 (45915) / 1.68
#  - the active cell ends up on B122 (the new "second to last" data cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the date that was already in row 122 (value only; format is unchanged).
$ws.Cells.Item(122, 1).Value = 45909

# Append the new observation in row 123, reusing the exact same cell
# formatting (date / 2-decimal number) as the row directly above it.
$ws.Range("A122:B122").Copy()
$ws.Range("A123:B123").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(123, 1).Value = 45915
$ws.Cells.Item(123, 2).Value = 1.68

# Mirror the cursor position recorded in the workbook: the active cell sits
# on B122 after the new row is appended below it.
$ws.Range("B122").Select()
